$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The photo-path column (B2/B3) originally held an absolute local path;
# it's now a relative path into a local "Fotos a enviar" folder.
$ws.Range("B2").Value = "Fotos a enviar\prueba.png"
$ws.Range("B3").Value = "Fotos a enviar\prueba.png"

# Leave the selection on B3, matching the saved view state.
[void]$ws.Range("B3").Select()
